# Scheduled-runner refresh: update cached market/profit figures on each
# sheet's Leve table (currentAveragePrice / LevePrice / LeveProfit columns)
# to the latest pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1414.4445
$ws.Range("J12").Value = 721.5
$ws.Range("L12").Value = 721.5
$ws.Range("N12").Value = -1061.5

$ws.Range("H15").Value = 2148.2036
$ws.Range("I15").Value = 2148.2036
$ws.Range("K15").Value = 6444.610799999999
$ws.Range("M15").Value = -6275.610799999999

$ws.Range("H33").Value = 275.5
$ws.Range("I33").Value = 240.53334
$ws.Range("K33").Value = 240.53334
$ws.Range("M33").Value = -11.53334000000001

$ws.Range("H40").Value = 1427.3334
$ws.Range("I40").Value = 1146.4445
$ws.Range("K40").Value = 1146.4445
$ws.Range("M40").Value = -971.4445000000001

$ws.Range("H41").Value = 242.60869
$ws.Range("I41").Value = 81
$ws.Range("J41").Value = 313.3125
$ws.Range("K41").Value = 81
$ws.Range("L41").Value = 313.3125
$ws.Range("M41").Value = 359
$ws.Range("N41").Value = -1193.3125

$ws.Range("H138").Value = 2503
$ws.Range("I138").Value = 1404.2273
$ws.Range("J138").Value = 3282.7742
$ws.Range("K138").Value = 4212.6819
$ws.Range("L138").Value = 9848.3226
$ws.Range("M138").Value = 927.3181000000004
$ws.Range("N138").Value = -20128.3226

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 166.66667
$ws.Range("I4").Value = 166.66667
$ws.Range("K4").Value = 166.66667
$ws.Range("M4").Value = -50.66667000000001

$ws.Range("H61").Value = 4673.5
$ws.Range("I61").Value = 3011.1428
$ws.Range("J61").Value = 5966.4443
$ws.Range("K61").Value = 3011.1428
$ws.Range("L61").Value = 5966.4443
$ws.Range("M61").Value = -2799.1428
$ws.Range("N61").Value = -6390.4443

$ws.Range("H108").Value = 32000
$ws.Range("J108").Value = 32000
$ws.Range("L108").Value = 32000
$ws.Range("N108").Value = -39680

$ws.Range("H112").Value = 34861.2
$ws.Range("J112").Value = 34861.2
$ws.Range("L112").Value = 34861.2
$ws.Range("N112").Value = -37815.2

$ws.Range("H122").Value = 1433.8636
$ws.Range("I122").Value = 1417.421
$ws.Range("K122").Value = 4252.263
$ws.Range("M122").Value = -1802.263

$ws.Range("H136").Value = 4673.5
$ws.Range("I136").Value = 3011.1428
$ws.Range("J136").Value = 5966.4443
$ws.Range("K136").Value = 9033.428400000001
$ws.Range("L136").Value = 17899.3329
$ws.Range("M136").Value = -6483.428400000001
$ws.Range("N136").Value = -22999.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 56311.895
$ws.Range("I134").Value = 81691.16
$ws.Range("J134").Value = 1323.5
$ws.Range("K134").Value = 245073.48
$ws.Range("L134").Value = 3970.5
$ws.Range("M134").Value = -242538.48
$ws.Range("N134").Value = -9040.5

$ws.Range("H139").Value = 44400
$ws.Range("I139").Value = 30000
$ws.Range("J139").Value = 48000
$ws.Range("K139").Value = 30000
$ws.Range("L139").Value = 48000
$ws.Range("M139").Value = -24860
$ws.Range("N139").Value = -58280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11401.475
$ws.Range("I31").Value = 18077.895
$ws.Range("J31").Value = 5360.905
$ws.Range("K31").Value = 18077.895
$ws.Range("L31").Value = 5360.905
$ws.Range("M31").Value = -17782.895
$ws.Range("N31").Value = -5950.905

$ws.Range("H34").Value = 11401.475
$ws.Range("I34").Value = 18077.895
$ws.Range("J34").Value = 5360.905
$ws.Range("K34").Value = 18077.895
$ws.Range("L34").Value = 5360.905
$ws.Range("M34").Value = -17875.895
$ws.Range("N34").Value = -5764.905

$ws.Range("H58").Value = 25138.38
$ws.Range("I58").Value = 1735.6
$ws.Range("J58").Value = 46413.637
$ws.Range("K58").Value = 1735.6
$ws.Range("L58").Value = 46413.637
$ws.Range("M58").Value = -1532.6
$ws.Range("N58").Value = -46819.637

$ws.Range("H136").Value = 25138.38
$ws.Range("I136").Value = 1735.6
$ws.Range("J136").Value = 46413.637
$ws.Range("K136").Value = 5206.799999999999
$ws.Range("L136").Value = 139240.911
$ws.Range("M136").Value = -2656.799999999999
$ws.Range("N136").Value = -144340.911

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2618
$ws.Range("I129").Value = 442.66666
$ws.Range("J129").Value = 4249.5
$ws.Range("K129").Value = 1327.99998
$ws.Range("L129").Value = 12748.5
$ws.Range("M129").Value = 3672.00002
$ws.Range("N129").Value = -22748.5

$ws.Range("H131").Value = 134179
$ws.Range("J131").Value = 139727.7
$ws.Range("L131").Value = 419183.1
$ws.Range("N131").Value = -429263.1

$ws.Range("H134").Value = 2488.8262
$ws.Range("I134").Value = 2145.3845
$ws.Range("J134").Value = 2935.3
$ws.Range("K134").Value = 6436.1535
$ws.Range("L134").Value = 8805.900000000001
$ws.Range("M134").Value = -1366.1535
$ws.Range("N134").Value = -18945.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3840.1177
$ws.Range("I126").Value = 3079.4119
$ws.Range("K126").Value = 9238.235700000001
$ws.Range("M126").Value = -6768.235700000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2785.72
$ws.Range("I7").Value = 2856.3333
$ws.Range("K7").Value = 2856.3333
$ws.Range("M7").Value = -2744.3333

$ws.Range("H46").Value = 1110.7059
$ws.Range("I46").Value = 799.38464
$ws.Range("K46").Value = 799.38464
$ws.Range("M46").Value = -611.38464

$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H126").Value = 2785.72
$ws.Range("I126").Value = 2856.3333
$ws.Range("K126").Value = 8568.999899999999
$ws.Range("M126").Value = -6098.999899999999

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1921.8077
$ws.Range("I132").Value = 1653.5
$ws.Range("J132").Value = 2525.5
$ws.Range("K132").Value = 4960.5
$ws.Range("L132").Value = 7576.5
$ws.Range("M132").Value = -2430.5
$ws.Range("N132").Value = -12636.5

$ws.Range("H136").Value = 1113585.9
$ws.Range("I136").Value = 1613812.2
$ws.Range("J136").Value = 1971.6666
$ws.Range("K136").Value = 4841436.6
$ws.Range("L136").Value = 5914.9998
$ws.Range("M136").Value = -4838886.6
$ws.Range("N136").Value = -11014.9998

$ws.Range("H139").Value = 52715
$ws.Range("J139").Value = 52715
$ws.Range("L139").Value = 52715
$ws.Range("N139").Value = -62995
